$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $text into $cellRef as a literal text value (never
# auto-converted to a number/date by Excel's type inference), without
# touching the cell's NumberFormat/style. Plain assignment via
# Range.Value coerces number-looking strings (e.g. "0.999", "1.00",
# "8.84") into real numbers, and forcing text via Range.NumberFormat
# ("@") leaves a permanent new cell style behind - neither of which
# matches the source data (plain inline/shared strings, default style).
#
# Instead we build a formula that evaluates to the literal string,
# write it into a scratch cell far outside the used range, copy that
# cell, and paste only the *value* (PasteSpecial xlPasteValues=-4163)
# into the target cell. Pasting a value (rather than assigning
# Range.Value directly) preserves the text type. The scratch cell is
# cleared afterwards so it leaves no trace in the saved workbook.
function Set-TextValue {
    param(
        [string]$cellRef,
        [string]$text
    )
    $escaped = $text.Replace('"', '""')
    $ws.Range("ZZ1").Formula = "=""$escaped"""
    $ws.Range("ZZ1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $ws.Range("ZZ1").ClearContents()
}

Set-TextValue "D2" "63.508.00"
Set-TextValue "E2" "  -2.78%  "
Set-TextValue "D3" "3.316.88"
Set-TextValue "E3" "  -4.54%  "
Set-TextValue "E4" "  +0.21%  "
Set-TextValue "D5" "547.98"
Set-TextValue "E5" "  -1.11%  "
Set-TextValue "D6" "171.56"
Set-TextValue "E6" "  -4.31%  "
Set-TextValue "E7" "  -4.00%  "
Set-TextValue "E8" "  +0.05%  "
Set-TextValue "E9" "  -3.89%  "
Set-TextValue "E10" "  -0.84%  "
Set-TextValue "E11" "  -1.87%  "
Set-TextValue "E12" "  -2.46%  "
Set-TextValue "D13" "8.84"
Set-TextValue "E13" "  -4.48%  "
Set-TextValue "D14" "3.846.16"
Set-TextValue "E14" "  -4.52%  "
Set-TextValue "E15" "  -3.34%  "
Set-TextValue "D16" "3.355.95"
Set-TextValue "E16" "  -3.17%  "
Set-TextValue "D17" "0.116"
Set-TextValue "E17" "  -3.77%  "
Set-TextValue "D18" "63.577.78"
Set-TextValue "E18" "  -2.56%  "
Set-TextValue "D19" "11.59"
Set-TextValue "E19" "  -3.25%  "
Set-TextValue "E20" "  -1.95%  "
Set-TextValue "D21" "411.53"
Set-TextValue "E21" "  -1.16%  "
Set-TextValue "D22" "4.03"
Set-TextValue "E22" "  -0.41%  "
Set-TextValue "D23" "4.38"
Set-TextValue "E23" "  +2.41%  "
Set-TextValue "D24" "13.70"
Set-TextValue "E24" "  +5.82%  "
Set-TextValue "D25" "82.72"
Set-TextValue "E25" "  -3.98%  "
Set-TextValue "E26" "  -3.25%  "
Set-TextValue "E27" "  -4.79%  "
Set-TextValue "D28" "8.58"
Set-TextValue "E28" "  -6.07%  "
Set-TextValue "D29" "28.96"
Set-TextValue "E29" "  -4.56%  "
Set-TextValue "E30" "  -3.43%  "
Set-TextValue "D31" "11.31"
Set-TextValue "E31" "  -4.03%  "
Set-TextValue "D32" "574.07"
Set-TextValue "E32" "  -5.56%  "
Set-TextValue "E33" "  -3.92%  "
Set-TextValue "D34" "57.50"
Set-TextValue "E34" "  -2.76%  "
Set-TextValue "E35" "  +0.11%  "
Set-TextValue "D36" "0.146"
Set-TextValue "E36" "  +0.86%  "
Set-TextValue "E37" "  -6.84%  "
Set-TextValue "D38" "3.40"
Set-TextValue "D39" "0.0₃0733"
Set-TextValue "E39" "  -7.13%  "
Set-TextValue "E40" "  -4.26%  "
Set-TextValue "D41" "3.113.80"
Set-TextValue "E41" "  -7.59%  "
Set-TextValue "D42" "0.999"
Set-TextValue "E42" "  +0.18%  "
Set-TextValue "B43" "ThetaToken"
Set-TextValue "C43" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D43" "2.76"
Set-TextValue "E43" "  -1.88%  "
Set-TextValue "B44" "ApeXProtocol"
Set-TextValue "C44" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D44" "3.23"
Set-TextValue "E44" "  -0.49%  "
Set-TextValue "D45" "0.0399"
Set-TextValue "E45" "  -3.54%  "
Set-TextValue "E46" "  -5.39%  "
Set-TextValue "D47" "2.59"
Set-TextValue "E47" "  -4.29%  "
Set-TextValue "E48" "  -4.01%  "
Set-TextValue "D49" "132.77"
Set-TextValue "E49" "  -3.66%  "
Set-TextValue "D50" "7.99"
Set-TextValue "E50" "  -5.54%  "
Set-TextValue "E51" "  +4.71%  "
